$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 258, shifting existing rows 258-353 down to 259-354
$ws.Rows.Item(258).Insert()

# Populate the newly inserted row 258 with its data
$ws.Cells.Item(258, 1).Value = 6
$ws.Cells.Item(258, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(258, 3).Value = "Metropolitana"
$ws.Cells.Item(258, 4).Value = 45146
$ws.Cells.Item(258, 5).Value = 13
$ws.Cells.Item(258, 6).Value = 100112029
$ws.Cells.Item(258, 7).Value = "Orégano"
$ws.Cells.Item(258, 8).Value = "Sin especificar"
$ws.Cells.Item(258, 9).Value = "Primera"
$ws.Cells.Item(258, 10).Value = 36
$ws.Cells.Item(258, 11).Value = 20000
$ws.Cells.Item(258, 12).Value = 20000
$ws.Cells.Item(258, 13).Value = 20000
$ws.Cells.Item(258, 14).Value = "$/docena de atados"
$ws.Cells.Item(258, 15).Value = "Región Metropolitana"
$ws.Cells.Item(258, 16).Value = 6667
$ws.Cells.Item(258, 17).Value = 3
$ws.Cells.Item(258, 18).Value = "Hortaliza"
